$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.999.08"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.789.27"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'359.69"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'109.76"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  -3.31%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "'40.22"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "'19.49"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "3.231.43"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "2.789.76"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "'0.942"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "51.940.83"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  -2.94%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("D23").Value = "'70.25"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'270.32"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E28").Value = "  +15.21%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "'0.0475"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'52.05"
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("D33").Value = "'34.37"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "  -5.09%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'19.07"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'119.74"
$ws.Range("E44").Value = "  -6.45%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.81"
$ws.Range("E45").Value = "  -9.00%  "
$ws.Range("D46").Value = "2.081.62"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("E47").Value = "  -4.78%  "
$ws.Range("D49").Value = "'5.78"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'0.957"
$ws.Range("E50").Value = "  -5.09%  "
$ws.Range("E51").Value = "  -2.93%  "
